$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Tighten the right indent on the paragraph that begins
#    "ด้วยข้าราชการครูและบุคลากรทางการศึกษา..."
#    rightChars 117 -> 58, right 281 (twips) -> 139 (twips)
# ---------------------------------------------------------------
$anchor = $d.Content.Duplicate
$null = $anchor.Find.Execute("ด้วยข้าราชการครูและบุคลากรทางการศึกษา", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $anchor.Paragraphs(1)
$targetPara.Format.RightIndent = 139 / 20.0
$targetPara.Format.CharacterUnitRightIndent = 58

# ---------------------------------------------------------------
# 2) Remove the space between "ราชการ" and "เพื่อ" in that sentence.
# ---------------------------------------------------------------
$null = $d.Content.Find.Execute( `
    "ด้วยข้าราชการครูและบุคลากรทางการศึกษาโรงเรียนวังน้ำเย็นวิทยาคมมีความประสงค์จะเดินทางไปราชการ เพื่อ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "ด้วยข้าราชการครูและบุคลากรทางการศึกษาโรงเรียนวังน้ำเย็นวิทยาคมมีความประสงค์จะเดินทางไปราชการเพื่อ", 2)

# ---------------------------------------------------------------
# 3) Re-add that space as its own run (matching formatting of the
#    following "{purpose} {date_range}" run, i.e. no complex-script
#    marker) immediately before "{purpose}". We do this by merging a
#    leading space into the "{purpose}" run via Find/Replace (which
#    inherits that run's own formatting) and then forcing the new
#    leading space back out into its own run by toggling a character
#    property on just that single character.
# ---------------------------------------------------------------
$null = $d.Content.Find.Execute("{purpose} {date_range}", $true, $false, $false, $false, $false, $true, 1, $false, " {purpose} {date_range}", 2)

$afterAnchor = $d.Content.Duplicate
$null = $afterAnchor.Find.Execute("{purpose}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newSpace = $d.Range($afterAnchor.Start - 1, $afterAnchor.Start)
$newSpace.Font.Bold = $true
$newSpace.Font.Bold = $false
